$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update the "Git" seminar cell to include the GitHub Classroom link
$ws.Range("E2").Value = "[Git](https://classroom.github.com/a/ojZuXTA9)"

# Reflect the new active cell selection as seen in the saved workbook
$ws.Range("E2").Select()
